# Regenerate merged AHB files
#
# 1) Header row (row 1) gets its "_old"/"_new" suffixes renamed to the
#    actual version tags being merged: "_FV2310" (left/base block,
#    columns A:J) and "_FV2404" (right/new block, columns L:U).
#    Column K ("diff") is unchanged.
# 2) The whole A1:U78 range becomes a named Excel Table ("Table1") with
#    its own autofilter.
# 3) The header row is frozen so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A..J (1..10): "<Name>_old"  -> "<Name>_FV2310"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2310"
}

# Columns L..U (12..21): "<Name>_new" -> "<Name>_FV2404"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2404"
}

# Freeze panes above row 2 (i.e. freeze the header row).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the full used range into an Excel Table with an autofilter.
# xlSrcRange = 1, XlYesNoGuess.xlYes = 1 (range already has headers).
$tableRange = $ws.Range("A1:U78")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
